$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The referentieregio (REFREG) list gains a new entry (REFREG14 gets split:
# the existing row keeps REFREG14 but with an updated provincie code, and a
# new row is inserted directly below it, also for REFREG14, carrying the
# value that used to belong to the original REFREG14 row). All subsequent
# rows (REFREG15, REFREG91, REFREG99) shift down by one row.

# Insert a new row at row 16, pushing REFREG15/REFREG91/REFREG99 down.
$ws.Rows("16").Insert()

# Row 15 (REFREG14) keeps its label but its value changes to 10000.
$ws.Range("B15").Value = 10000

# The newly inserted row 16 becomes a second REFREG14 entry with the value
# that row 15 used to have (40000).
$ws.Range("A16").Value = "REFREG14"
$ws.Range("B16").Value = 40000
